$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 9: 10-Aug-2023 attendance entry
$ws.Range("A9").Value = 45148
$ws.Range("A9").NumberFormat = "d-mmm-yy"

$ws.Range("B9").Value = "PRESENT"
$ws.Range("C9").Value = "PRESENT"
$ws.Range("D9").Value = "PRESENT"
$ws.Range("E9").Value = "PRESENT"
$ws.Range("F9").Value = "PRESENT"
$ws.Range("G9").Value = "PRESENT"
$ws.Range("H9").Value = "PRESENT"
$ws.Range("I9").Value = "PRESENT"
$ws.Range("J9").Value = "PRESENT"
$ws.Range("K9").Value = "ABSENT"

# Comment explaining the absence on K9
$excel.UserName = "LENOVO"
$comment = $ws.Range("K9").AddComment("LENOVO:" + [char]10 + "Document work")

# Move the selection the way the saved workbook shows it
$ws.Range("B9").Select() | Out-Null
